$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a far-away scratch cell to stage each new value as forced Text
# (NumberFormat "@" prevents Excel from auto-coercing numeric-looking
# strings like "0.9998" or "1.000" into floating point numbers, which
# would corrupt exact text such as trailing zeros). We then copy only
# the *value* (PasteSpecial xlPasteValues = -4163) into the target cell
# so the target keeps its original (unstyled) formatting, matching the
# workbook before this edit.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"

$scratch.Value = '29.395.07'
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)

$scratch.Value = '1.847.74'
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)

$scratch.Value = '  +0.23%  '
$scratch.Copy()
$ws.Range("E3").PasteSpecial(-4163)

$scratch.Value = '0.9998'
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)

$scratch.Value = '  +0.09%  '
$scratch.Copy()
$ws.Range("E4").PasteSpecial(-4163)

$scratch.Value = '240.38'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)

$scratch.Value = '  +0.08%  '
$scratch.Copy()
$ws.Range("E5").PasteSpecial(-4163)

$scratch.Value = '0.6295'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)

$scratch.Value = '  -0.70%  '
$scratch.Copy()
$ws.Range("E6").PasteSpecial(-4163)

$scratch.Value = '  +0.07%  '
$scratch.Copy()
$ws.Range("E7").PasteSpecial(-4163)

$scratch.Value = '0.07606'
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)

$scratch.Value = '  +0.87%  '
$scratch.Copy()
$ws.Range("E8").PasteSpecial(-4163)

$scratch.Value = '0.2930'
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)

$scratch.Value = '  -0.95%  '
$scratch.Copy()
$ws.Range("E9").PasteSpecial(-4163)

$scratch.Value = '  -0.74%  '
$scratch.Copy()
$ws.Range("E10").PasteSpecial(-4163)

$scratch.Value = '0.07738'
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)

$scratch.Value = '  +0.10%  '
$scratch.Copy()
$ws.Range("E11").PasteSpecial(-4163)

$scratch.Value = '1.849.43'
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)

$scratch.Value = '  -6.83%  '
$scratch.Copy()
$ws.Range("E12").PasteSpecial(-4163)

$scratch.Value = '5.002'
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)

$scratch.Value = '  +0.35%  '
$scratch.Copy()
$ws.Range("E13").PasteSpecial(-4163)

$scratch.Value = '0.00001086'
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)

$scratch.Value = '  +9.77%  '
$scratch.Copy()
$ws.Range("E14").PasteSpecial(-4163)

$scratch.Value = '0.6789'
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)

$scratch.Value = '  -0.57%  '
$scratch.Copy()
$ws.Range("E15").PasteSpecial(-4163)

$scratch.Value = '  +0.67%  '
$scratch.Copy()
$ws.Range("E16").PasteSpecial(-4163)

$scratch.Value = '2.100.64'
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)

$scratch.Value = '  -7.23%  '
$scratch.Copy()
$ws.Range("E17").PasteSpecial(-4163)

$scratch.Value = '6.174'
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)

$scratch.Value = '  +0.19%  '
$scratch.Copy()
$ws.Range("E18").PasteSpecial(-4163)

$scratch.Value = '29.412.39'
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)

$scratch.Value = '  +0.06%  '
$scratch.Copy()
$ws.Range("E19").PasteSpecial(-4163)

$scratch.Value = '228.49'
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)

$scratch.Value = '  -0.82%  '
$scratch.Copy()
$ws.Range("E20").PasteSpecial(-4163)

$scratch.Value = '  +0.02%  '
$scratch.Copy()
$ws.Range("E21").PasteSpecial(-4163)

$scratch.Value = '  +0.09%  '
$scratch.Copy()
$ws.Range("E22").PasteSpecial(-4163)

$scratch.Value = '7.479'
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)

$scratch.Value = '  -0.80%  '
$scratch.Copy()
$ws.Range("E23").PasteSpecial(-4163)

$scratch.Value = '1.001'
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)

$scratch.Value = '  +0.07%  '
$scratch.Copy()
$ws.Range("E24").PasteSpecial(-4163)

$scratch.Value = '157.31'
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)

$scratch.Value = '  +0.78%  '
$scratch.Copy()
$ws.Range("E25").PasteSpecial(-4163)

$scratch.Value = '0.1397'
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)

$scratch.Value = '  -0.65%  '
$scratch.Copy()
$ws.Range("E26").PasteSpecial(-4163)

$scratch.Value = '8.344'
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)

$scratch.Value = '17.63'
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)

$scratch.Value = '  -0.17%  '
$scratch.Copy()
$ws.Range("E28").PasteSpecial(-4163)

$scratch.Value = '1.466'
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)

$scratch.Value = '  -0.20%  '
$scratch.Copy()
$ws.Range("E29").PasteSpecial(-4163)

$scratch.Value = '0.05585'
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)

$scratch.Value = '  -1.95%  '
$scratch.Copy()
$ws.Range("E31").PasteSpecial(-4163)

$scratch.Value = '4.105'
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)

$scratch.Value = '  -0.40%  '
$scratch.Copy()
$ws.Range("E32").PasteSpecial(-4163)

$scratch.Value = '4.036'
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)

$scratch.Value = '  +0.25%  '
$scratch.Copy()
$ws.Range("E33").PasteSpecial(-4163)

$scratch.Value = '1.844'
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)

$scratch.Value = '  +0.03%  '
$scratch.Copy()
$ws.Range("E34").PasteSpecial(-4163)

$scratch.Value = '  +0.06%  '
$scratch.Copy()
$ws.Range("E35").PasteSpecial(-4163)

$scratch.Value = '0.7097'
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)

$scratch.Value = '  -0.67%  '
$scratch.Copy()
$ws.Range("E36").PasteSpecial(-4163)

$scratch.Value = '2.587'
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)

$scratch.Value = '  -0.25%  '
$scratch.Copy()
$ws.Range("E37").PasteSpecial(-4163)

$scratch.Value = '1.239.97'
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)

$scratch.Value = '  -0.81%  '
$scratch.Copy()
$ws.Range("E38").PasteSpecial(-4163)

$scratch.Value = '2.775'
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)

$scratch.Value = '  -0.93%  '
$scratch.Copy()
$ws.Range("E39").PasteSpecial(-4163)

$scratch.Value = '0.01802'
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)

$scratch.Value = '  -0.54%  '
$scratch.Copy()
$ws.Range("E40").PasteSpecial(-4163)

$scratch.Value = '6.413'
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)

$scratch.Value = '  +5.26%  '
$scratch.Copy()
$ws.Range("E41").PasteSpecial(-4163)

$scratch.Value = '0.9047'
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)

$scratch.Value = '  +0.32%  '
$scratch.Copy()
$ws.Range("E42").PasteSpecial(-4163)

$scratch.Value = '1.000'
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)

$scratch.Value = '  +0.09%  '
$scratch.Copy()
$ws.Range("E43").PasteSpecial(-4163)

$scratch.Value = '101.81'
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)

$scratch.Value = '  -0.06%  '
$scratch.Copy()
$ws.Range("E44").PasteSpecial(-4163)

$scratch.Value = '65.98'
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)

$scratch.Value = '  -0.50%  '
$scratch.Copy()
$ws.Range("E45").PasteSpecial(-4163)

$scratch.Value = '  +3.59%  '
$scratch.Copy()
$ws.Range("E46").PasteSpecial(-4163)

$scratch.Value = '7.176'
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)

$scratch.Value = '  +1.64%  '
$scratch.Copy()
$ws.Range("E47").PasteSpecial(-4163)

$scratch.Value = '0.4022'
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)

$scratch.Value = '  +0.14%  '
$scratch.Copy()
$ws.Range("E48").PasteSpecial(-4163)

$scratch.Value = '9.024'
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)

$scratch.Value = '  -1.40%  '
$scratch.Copy()
$ws.Range("E49").PasteSpecial(-4163)

$scratch.Value = '1.678'
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)

$scratch.Value = '  -1.42%  '
$scratch.Copy()
$ws.Range("E50").PasteSpecial(-4163)

$scratch.Value = '0.1120'
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)

$scratch.Value = '  -0.40%  '
$scratch.Copy()
$ws.Range("E51").PasteSpecial(-4163)

$scratch.Clear()
$excel.CutCopyMode = 0
